# edit.ps1
# Reproduces the commit:
#   1. The table on slide 5 switches from the deck's local/default table
#      style ({7AEC0FE1-...}) to the built-in PowerPoint table style
#      {363D1559-32FE-4298-BFEA-2F2F459DD4EF} ("Medium Style 2 - Accent 1").
#   2. The deck's main theme (ppt/theme/theme1.xml, the one used by the
#      Slide Master / all slides) swaps its 12-slot colour scheme from the
#      "Red Violet" palette ("Integral" theme) to the stock "Office"
#      palette ("Office Theme"). (The font scheme and format scheme were
#      already identical between the two themes in this deck, so only the
#      colours actually change.)

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ------------------------------------------

$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{363D1559-32FE-4298-BFEA-2F2F459DD4EF}")
    }
}

# --- 2. Theme colour scheme swap (Red Violet -> Office) ------------------

# The presentation has a single Slide Master / theme shared by every
# slide, so editing it through any slide's ThemeColorScheme updates the
# underlying theme part (ppt/theme/theme1.xml) for the whole deck.
$themeSlide = $p.Slides.Item(1)
$colors = $themeSlide.ThemeColorScheme

# Index order follows MsoThemeColorSchemeIndex:
#  1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#  5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink
$officeColors = @(
    0x000000,  # Dark1
    0xFFFFFF,  # Light1
    0x44546A,  # Dark2
    0xE7E6E6,  # Light2
    0x5B9BD5,  # Accent1
    0xED7D31,  # Accent2
    0xA5A5A5,  # Accent3
    0xFFC000,  # Accent4
    0x4472C4,  # Accent5
    0x70AD47,  # Accent6
    0x0563C1,  # Hyperlink
    0x954F72   # FollowedHyperlink
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $rgbHex = $officeColors[$i - 1]
    # COM RGB() is stored little-endian as 0x00BBGGRR
    $r = ($rgbHex -band 0xFF0000) -shr 16
    $g = ($rgbHex -band 0x00FF00) -shr 8
    $b = ($rgbHex -band 0x0000FF)
    $comRgb = $b * 65536 + $g * 256 + $r
    $colors.Item($i).RGB = $comRgb
}
